# Update the "Analysis Results" battery-dashboard sheet:
#  - rows 6-7: swap Starting/Ending SoC(%) values
#  - rows 8-33: append unit suffixes to several labels, relabel/reorder a
#    few rows (cell-voltage & temperature pairs swapped, BMS-temp row
#    replaced by Battery Voltage, etc.) and refresh the numbers that moved
#    along with them
#  - rows 34-42: the speed-bucket ("Time spent in X-Y km/h") block shifts
#    down by one row with refreshed percentages
#  - row 43: new "Time spent in 80-90 km/h" row added, extending the used
#    range to A1:B43
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'Starting SoC (%)'
$ws.Range("B6").Value = 99

$ws.Range("A7").Value = 'Ending SoC (%)'
$ws.Range("B7").Value = 17

$ws.Range("A8").Value = 'Total distance covered (km)'
$ws.Range("B8").Value = 36.21829742812135

$ws.Range("A9").Value = 'Total energy consumption(WH/KM)'
$ws.Range("B9").Value = 45.34236056445818

$ws.Range("A10").Value = 'Total SOC consumed(%)'
$ws.Range("B10").Value = 82

$ws.Range("A12").Value = 'Peak Power(kW)'
$ws.Range("B12").Value = 5233.3208

$ws.Range("A13").Value = 'Average Power(kW)'
$ws.Range("B13").Value = -1489.920152132056

$ws.Range("A14").Value = 'Total Energy Regenerated(kWh)'
$ws.Range("B14").Value = 3.365682143888889

$ws.Range("A15").Value = 'Regenerative Effectiveness(%)'
$ws.Range("B15").Value = 0.2045275331437714

$ws.Range("A16").Value = 'Highest Cell Voltage(V)'
$ws.Range("B16").Value = 3.491

$ws.Range("A17").Value = 'Lowest Cell Voltage(V)'
$ws.Range("B17").Value = 3.032

$ws.Range("A18").Value = 'Difference in Cell Voltage(V)'
$ws.Range("B18").Value = 0.4590000000000001

$ws.Range("A19").Value = 'Minimum Temperature(C)'
$ws.Range("B19").Value = 31

$ws.Range("A20").Value = 'Maximum Temperature(C)'
$ws.Range("B20").Value = 47

$ws.Range("A21").Value = 'Difference in Temperature(C)'
$ws.Range("B21").Value = 16

$ws.Range("A22").Value = 'Maximum Fet Temperature-BMS(C)'
$ws.Range("B22").Value = 60

$ws.Range("A23").Value = 'Maximum Afe Temperature-BMS(C)'
$ws.Range("B23").Value = 63

$ws.Range("A24").Value = 'Maximum PCB Temperature-BMS(C)'
$ws.Range("B24").Value = 60

$ws.Range("A25").Value = 'Maximum MCU Temperature(C)'
$ws.Range("B25").Value = 59

$ws.Range("A26").Value = 'Maximum Motor Temperature(C)'
$ws.Range("B26").Value = 98

$ws.Range("A27").Value = 'Abnormal Motor Temperature Detected(C)'
$ws.Range("B27").Value = 0

$ws.Range("A28").Value = 'highest cell temp(C)'
$ws.Range("B28").Value = 47

$ws.Range("A29").Value = 'lowest cell temp(C)'
$ws.Range("B29").Value = 31

$ws.Range("A30").Value = 'Difference between Highest and Lowest Cell Temperature at 100% SOC(C)'
$ws.Range("B30").Value = 16

$ws.Range("A31").Value = 'Battery Voltage(V)'
$ws.Range("B31").Value = 54

$ws.Range("A32").Value = 'Total energy charged(kWh)'
$ws.Range("B32").Value = 1.746260955

$ws.Range("A33").Value = 'Electricity consumption units(kW)'
$ws.Range("B33").Value = 0.0000001215720520050125

$ws.Range("A34").Value = 'Idling time percentage'
$ws.Range("B34").Value = 24.73596594165325

$ws.Range("A35").Value = 'Time spent in 0-10 km/h'
$ws.Range("B35").Value = 5.108752012662719

$ws.Range("A36").Value = 'Time spent in 10-20 km/h'
$ws.Range("B36").Value = 6.153972109270529

$ws.Range("A37").Value = 'Time spent in 20-30 km/h'
$ws.Range("B37").Value = 9.349671151379527

$ws.Range("A38").Value = 'Time spent in 30-40 km/h'
$ws.Range("B38").Value = 14.190977812952

$ws.Range("A39").Value = 'Time spent in 40-50 km/h'
$ws.Range("B39").Value = 8.708348115601888

$ws.Range("A40").Value = 'Time spent in 50-60 km/h'
$ws.Range("B40").Value = 6.536036896542313

$ws.Range("A41").Value = 'Time spent in 60-70 km/h'
$ws.Range("B41").Value = 14.48571350599023

$ws.Range("A42").Value = 'Time spent in 70-80 km/h'
$ws.Range("B42").Value = 10.61867205196081

$ws.Range("A43").Value = 'Time spent in 80-90 km/h'
$ws.Range("B43").Value = 0
